# Auto-generated edit script applying the cryptos.xlsx diff (cryptocurrency price/volume update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.743.73'
$ws.Range("E2").Value = '  -5.19%  '
$ws.Range("D3").Value = '3.363.21'
$ws.Range("E3").Value = '  -6.59%  '
$ws.Range("E4").Value = '  +0.06%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '559.93'
$c.ClearFormats()
$ws.Range("E5").Value = '  -5.70%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '183.69'
$c.ClearFormats()
$ws.Range("E6").Value = '  -8.51%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.599'
$c.ClearFormats()
$ws.Range("E7").Value = '  -4.53%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '3.353.62'
$ws.Range("E9").Value = '  -6.51%  '
$ws.Range("E10").Value = '  -12.91%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.595'
$c.ClearFormats()
$ws.Range("E11").Value = '  -7.74%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '48.00'
$c.ClearFormats()
$ws.Range("E12").Value = '  -10.23%  '
$ws.Range("E13").Value = '  -10.29%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '8.74'
$c.ClearFormats()
$ws.Range("E14").Value = '  -9.75%  '
$ws.Range("D15").Value = '3.899.86'
$ws.Range("E15").Value = '  -6.52%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '603.70'
$c.ClearFormats()
$ws.Range("E16").Value = '  -12.96%  '
$ws.Range("D17").Value = '66.554.70'
$ws.Range("E17").Value = '  -5.58%  '
$ws.Range("D18").Value = '3.366.43'
$ws.Range("E18").Value = '  -6.65%  '
$ws.Range("E19").Value = '  -4.50%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '17.75'
$c.ClearFormats()
$ws.Range("E20").Value = '  -6.93%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '11.67'
$c.ClearFormats()
$ws.Range("E21").Value = '  -8.48%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.914'
$c.ClearFormats()
$ws.Range("E22").Value = '  -8.17%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '16.98'
$c.ClearFormats()
$ws.Range("E23").Value = '  -7.59%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '5.08'
$c.ClearFormats()
$ws.Range("E24").Value = '  -3.86%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '96.97'
$c.ClearFormats()
$ws.Range("E25").Value = '  -12.35%  '
$ws.Range("E26").Value = '  -10.19%  '
$ws.Range("E27").Value = '  -8.70%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '9.54'
$c.ClearFormats()
$ws.Range("E28").Value = '  -9.41%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '8.83'
$c.ClearFormats()
$ws.Range("E29").Value = '  -12.08%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '30.78'
$c.ClearFormats()
$ws.Range("E30").Value = '  -11.08%  '
$ws.Range("E31").Value = '  -12.83%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '6.36'
$c.ClearFormats()
$ws.Range("E32").Value = '  -10.13%  '
$ws.Range("E33").Value = '  -8.83%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.106'
$c.ClearFormats()
$ws.Range("E34").Value = '  -7.44%  '
$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '545.72'
$c.ClearFormats()
$ws.Range("E35").Value = '  +6.99%  '
$ws.Range("D36").Value = '3.821.80'
$ws.Range("E36").Value = '  -0.14%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '58.41'
$c.ClearFormats()
$ws.Range("E37").Value = '  -8.28%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range("E38").Value = '  +0.04%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '3.76'
$c.ClearFormats()
$ws.Range("E39").Value = '  +39.04%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '3.40'
$c.ClearFormats()
$ws.Range("E40").Value = '  -6.42%  '
$ws.Range("E41").Value = '  -14.61%  '
$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.73'
$c.ClearFormats()
$ws.Range("E42").Value = '  -9.54%  '
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.129'
$c.ClearFormats()
$ws.Range("E43").Value = '  -6.12%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.352'
$c.ClearFormats()
$ws.Range("E44").Value = '  -7.91%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '32.70'
$c.ClearFormats()
$ws.Range("E45").Value = '  -10.86%  '
$ws.Range("E46").Value = '  -11.45%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '2.69'
$c.ClearFormats()
$ws.Range("E47").Value = '  -12.34%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '3.13'
$c.ClearFormats()
$ws.Range("E48").Value = '  -9.43%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.131'
$c.ClearFormats()
$ws.Range("E49").Value = '  -7.23%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.ClearFormats()
$ws.Range("E50").Value = '  -0.35%  '
$ws.Range("E51").Value = '  -10.49%  '
